# Update "想去人数" (F column) counts on the "展览", "演出" and "全部类型"
# sheets to reflect newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 10197   # was 10178
$wsExpo.Range("F4").Value = 54      # was 52
$wsExpo.Range("F5").Value = 631     # was 628
$wsExpo.Range("F6").Value = 482     # was 480

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 5       # was 4

# --- Sheet "全部类型" (All types, aggregated) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10197    # was 10178
$wsAll.Range("F4").Value = 54       # was 52
$wsAll.Range("F5").Value = 631      # was 628
$wsAll.Range("F6").Value = 5        # was 4
$wsAll.Range("F7").Value = 482      # was 480
